$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.474.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.477.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +3.34%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.927.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.306.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.468.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0908"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "504.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.328"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0258"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.584"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
